$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) : new columns E (Blog URL) and F (Blog RSS) ----
# Match the existing header look (font size 18, same as A1:D1 -> style index 2)
$ws.Range("E1").Value = "Blog URL"
$ws.Range("E1").Font.Size = 18
$ws.Range("F1").Value = "Blog RSS"
$ws.Range("F1").Font.Size = 18

# ---- Data rows: Blog URL (col E, hyperlinked) + Blog RSS (col F, plain text) ----
# Row 2 - Ben Dudley
$ws.Range("F2").Value = "https://bed19.wordpress.com/feed/"

# Row 3 - David Fairbrother
$ws.Range("F3").Value = "http://users.aber.ac.uk/daf5/blog/?feed=rss2"

# Row 4 - Jonathan Englund
$ws.Range("F4").Value = "http://users.aber.ac.uk/jee17/wordpress/?feed=rss2"

# Row 5 - Josh Doyle
$ws.Range("F5").Value = "http://jod32.blogspot.com/feeds/posts/default"

# Row 6 - Liam Fitzgerald
$ws.Range("F6").Value = "http://users.aber.ac.uk/lif5/wordpress/?feed=rss2"

# Row 7 - Maurice Corriette (no blog info)

# Row 8 - Olver Earl
$ws.Range("F8").Value = "http://users.aber.ac.uk/ole4/group-project/?feed=rss2"

# Row 9 - Tim Anderson
$ws.Range("F9").Value = "http://users.aber.ac.uk/tma1/wordpress/?feed=rss2"

# ---- Hyperlinks for column E (this also writes the display text + shared strings) ----
$ws.Hyperlinks.Add($ws.Range("E8"), "http://users.aber.ac.uk/ole4/group-project")
$ws.Hyperlinks.Add($ws.Range("E6"), "http://users.aber.ac.uk/lif5/wordpress/")
$ws.Hyperlinks.Add($ws.Range("E9"), "http://users.aber.ac.uk/tma1/wordpress/")
$ws.Hyperlinks.Add($ws.Range("E4"), "http://users.aber.ac.uk/jee17/wordpress/")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://bed19.wordpress.com/")
$ws.Hyperlinks.Add($ws.Range("E5"), "http://jod32.blogspot.co.uk/")
$ws.Hyperlinks.Add($ws.Range("E3"), "http://users.aber.ac.uk/daf5/blog/")

# ---- Column widths (bestFit-like, matching the author's manual column resize) ----
$ws.Columns.Item(5).ColumnWidth = 39.3
$ws.Columns.Item(6).ColumnWidth = 48.15

# ---- Selection moved to F3 (where the author ended up) ----
$ws.Range("F3").Select()
